$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 266.66666
$ws.Range("I9").Value = 307.5
$ws.Range("K9").Value = 307.5
$ws.Range("M9").Value = -138.5
$ws.Range("H46").Value = 55559532
$ws.Range("J46").Value = 5959.5
$ws.Range("L46").Value = 17878.5
$ws.Range("N46").Value = -18116.5
$ws.Range("H48").Value = 4615.3076
$ws.Range("J48").Value = 4615.3076
$ws.Range("L48").Value = 13845.9228
$ws.Range("N48").Value = -14429.9228
$ws.Range("H56").Value = 4615.3076
$ws.Range("J56").Value = 4615.3076
$ws.Range("L56").Value = 13845.9228
$ws.Range("N56").Value = -14913.9228
$ws.Range("H60").Value = 55559532
$ws.Range("J60").Value = 5959.5
$ws.Range("L60").Value = 17878.5
$ws.Range("N60").Value = -18846.5
$ws.Range("H74").Value = 4317.3125
$ws.Range("I74").Value = 3957.7778
$ws.Range("J74").Value = 4779.5713
$ws.Range("K74").Value = 3957.7778
$ws.Range("L74").Value = 4779.5713
$ws.Range("M74").Value = -3021.7778
$ws.Range("N74").Value = -6651.5713
$ws.Range("H77").Value = 4317.3125
$ws.Range("I77").Value = 3957.7778
$ws.Range("J77").Value = 4779.5713
$ws.Range("K77").Value = 19788.889
$ws.Range("L77").Value = 23897.8565
$ws.Range("M77").Value = -15108.889
$ws.Range("N77").Value = -33257.85649999999
$ws.Range("H137").Value = 14707975
$ws.Range("I137").Value = 21741364
$ws.Range("J137").Value = 1795.7273
$ws.Range("K137").Value = 65224092
$ws.Range("L137").Value = 5387.1819
$ws.Range("M137").Value = -65221542
$ws.Range("N137").Value = -10487.1819
$ws.Range("H138").Value = 1265.1143
$ws.Range("I138").Value = 982.63336
$ws.Range("J138").Value = 2960
$ws.Range("K138").Value = 2947.90008
$ws.Range("L138").Value = 8880
$ws.Range("M138").Value = 2192.09992
$ws.Range("N138").Value = -19160

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5963.209
$ws.Range("I32").Value = 6326.45
$ws.Range("J32").Value = 2849.7144
$ws.Range("K32").Value = 6326.45
$ws.Range("L32").Value = 2849.7144
$ws.Range("M32").Value = -6039.45
$ws.Range("N32").Value = -3423.7144
$ws.Range("H81").Value = 23590.5
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 23590.5
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 23590.5
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -25586.5
$ws.Range("H84").Value = 23590.5
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 23590.5
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 70771.5
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -80755.5

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3153.5625
$ws.Range("I20").Value = 2266.375
$ws.Range("J20").Value = 4040.75
$ws.Range("K20").Value = 2266.375
$ws.Range("L20").Value = 4040.75
$ws.Range("M20").Value = -2019.375
$ws.Range("N20").Value = -4534.75
$ws.Range("H134").Value = 4306.5967
$ws.Range("I134").Value = 1923.6216
$ws.Range("J134").Value = 8715.1
$ws.Range("K134").Value = 5770.864799999999
$ws.Range("L134").Value = 26145.3
$ws.Range("M134").Value = -3235.864799999999
$ws.Range("N134").Value = -31215.3

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 8000
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1021.2941
$ws.Range("I5").Value = 290.16666
$ws.Range("J5").Value = 2776
$ws.Range("K5").Value = 870.4999799999999
$ws.Range("L5").Value = 8328
$ws.Range("M5").Value = -758.4999799999999
$ws.Range("N5").Value = -8552
$ws.Range("H47").Value = 2748.375
$ws.Range("I47").Value = 180.6
$ws.Range("J47").Value = 7028
$ws.Range("K47").Value = 541.8
$ws.Range("L47").Value = 21084
$ws.Range("M47").Value = -110.8
$ws.Range("N47").Value = -21946
$ws.Range("H51").Value = 1662.5
$ws.Range("I51").Value = 757.1429000000001
$ws.Range("K51").Value = 2271.4287
$ws.Range("M51").Value = -1811.4287
$ws.Range("H118").Value = 4684.875
$ws.Range("I118").Value = 619.75
$ws.Range("J118").Value = 8750
$ws.Range("K118").Value = 1859.25
$ws.Range("L118").Value = 26250
$ws.Range("M118").Value = -616.25
$ws.Range("N118").Value = -28736
$ws.Range("H131").Value = 970.42426
$ws.Range("I131").Value = 626.25
$ws.Range("J131").Value = 1080.56
$ws.Range("K131").Value = 1878.75
$ws.Range("L131").Value = 3241.68
$ws.Range("M131").Value = 3161.25
$ws.Range("N131").Value = -13321.68
$ws.Range("H135").Value = 1021.2941
$ws.Range("I135").Value = 290.16666
$ws.Range("J135").Value = 2776
$ws.Range("K135").Value = 2611.49994
$ws.Range("L135").Value = 24984
$ws.Range("M135").Value = -76.4999399999997
$ws.Range("N135").Value = -30054

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2712.5
$ws.Range("I80").Value = 2616.6667
$ws.Range("K80").Value = 2616.6667
$ws.Range("M80").Value = -1618.6667
$ws.Range("H83").Value = 2712.5
$ws.Range("I83").Value = 2616.6667
$ws.Range("K83").Value = 13083.3335
$ws.Range("M83").Value = -8091.333500000001

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 675.3158
$ws.Range("J22").Value = 654.35
$ws.Range("L22").Value = 654.35
$ws.Range("N22").Value = -1244.35
$ws.Range("H27").Value = 675.3158
$ws.Range("J27").Value = 654.35
$ws.Range("L27").Value = 654.35
$ws.Range("N27").Value = -868.35
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H74").Value = 18980
$ws.Range("I74").Value = 18900
$ws.Range("J74").Value = 19000
$ws.Range("K74").Value = 18900
$ws.Range("L74").Value = 19000
$ws.Range("M74").Value = -17902
$ws.Range("N74").Value = -20996
$ws.Range("H77").Value = 18980
$ws.Range("I77").Value = 18900
$ws.Range("J77").Value = 19000
$ws.Range("K77").Value = 56700
$ws.Range("L77").Value = 57000
$ws.Range("M77").Value = -51708
$ws.Range("N77").Value = -66984
$ws.Range("H82").Value = 2414
$ws.Range("I82").Value = 2184.875
$ws.Range("J82").Value = 2780.6
$ws.Range("K82").Value = 2184.875
$ws.Range("L82").Value = 2780.6
$ws.Range("M82").Value = -1823.875
$ws.Range("N82").Value = -3502.6
$ws.Range("H85").Value = 2414
$ws.Range("I85").Value = 2184.875
$ws.Range("J85").Value = 2780.6
$ws.Range("K85").Value = 2184.875
$ws.Range("L85").Value = 2780.6
$ws.Range("M85").Value = -936.875
$ws.Range("N85").Value = -5276.6

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 6000
$ws.Range("I61").Value = 6000
$ws.Range("K61").Value = 6000
$ws.Range("M61").Value = -5664
